# Rename the existing (only) worksheet from "test_youtuber" to "Template".
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Template"

# Add a new worksheet "Available Fields" right after "Template".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Available Fields"

# Put the wiki-link text in A1 and turn it into a real hyperlink.
$ws2.Range("A1").Value = "See https://github.com/71tech/ODIN/wiki/Available-fields"
$ws2.Hyperlinks.Add($ws2.Range("A1"), "https://github.com/71tech/ODIN/wiki/Available-fields") | Out-Null
